$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The PlannedQty column (G) stores numeric-looking values ("665", "661",
# "685") as TEXT (shared strings), not numbers. A plain .Value = "685"
# assignment would have Excel re-infer it as a genuine number, so instead
# we copy the already-text-typed value from an untouched reference row
# (G2/G3/G4 are outside the edited range) and paste just the value -
# this preserves both the destination cell's existing style and the
# source's text typing.
$xlPasteValues = -4163

# Row 6: becomes the "Press Approval Task" block (previously held by row 8)
$ws.Range("B6").Value = "-"
$ws.Range("D6").Value = "169-Press Approval Task "
$ws.Range("G4").Copy()
$ws.Range("G6").PasteSpecial($xlPasteValues)
$ws.Range("L6").Value = "Press Approval Task"
$ws.Range("M6").Value = "169-Press Approval Task "

# Row 7: becomes the "Digital Print" block (previously held by row 6)
$ws.Range("B7").Value = "Digital Print F 4x0"
$ws.Range("D7").Value = "252-HP 10000 Press"
$ws.Range("G2").Copy()
$ws.Range("G7").PasteSpecial($xlPasteValues)
$ws.Range("L7").Value = "252-HP 10000 Press"
$ws.Range("M7").Value = "252-HP 10000 Press"

# Row 8: becomes the "Cut" block (previously held by row 7)
$ws.Range("B8").Value = "Cut"
$ws.Range("D8").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("G3").Copy()
$ws.Range("G8").PasteSpecial($xlPasteValues)
$ws.Range("L8").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("M8").Value = "406-45`" Polar 115ED Cutter`n404-45`" Polar 115EMC Cutter`n405-54`" Polar 137EMC Cutter`n402-45`" Polar 115EMC Cutter`n403-54`" Polar 137ED Cutter"

$excel.CutCopyMode = $false

